{"js": "// Replace the 25 two-digit-division answers in the 5x5 \"data\" rows of the\n// worksheet table (rows 0, 4, 8, 12, 16 of the 20-row table; the table has\n// 4 blank rows between each data row). Each cell's text is replaced in\n// document order with its corresponding new value from the commit.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map: table row index -> new values for that row's 5 columns, in order.\nconst newRowValues = {\n  0: [\"51\u00f74=12, 3\", \"20\u00f73=6, 2\", \"15\u00f73=5, 0\", \"10\u00f78=1, 2\", \"88\u00f78=11, 0\"],\n  4: [\"34\u00f79=3, 7\", \"14\u00f78=1, 6\", \"48\u00f74=12, 0\", \"40\u00f73=13, 1\", \"31\u00f78=3, 7\"],\n  8: [\"98\u00f73=32, 2\", \"30\u00f74=7, 2\", \"31\u00f76=5, 1\", \"97\u00f79=10, 7\", \"13\u00f73=4, 1\"],\n  12: [\"82\u00f74=20, 2\", \"24\u00f74=6, 0\", \"64\u00f72=32, 0\", \"33\u00f78=4, 1\", \"88\u00f75=17, 3\"],\n  16: [\"72\u00f74=18, 0\", \"15\u00f73=5, 0\", \"67\u00f76=11, 1\", \"37\u00f77=5, 2\", \"60\u00f72=30, 0\"],\n};\n\nfor (const rowIndexStr of Object.keys(newRowValues)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const rowValues = newRowValues[rowIndex];\n  for (let col = 0; col < rowValues.length; col++) {\n    table.getCell(rowIndex, col).value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 two-digit-division answers in the 5x5 \"data\" rows of the\n# worksheet table (rows 1, 5, 9, 13, 17 of the 20-row table, 1-based; the\n# table has 4 blank rows between each data row). Each cell's text is\n# replaced in document order with its corresponding new value.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$dataRows = @(1, 5, 9, 13, 17)\n$newValues = @{\n    1  = @(\"51\u00f74=12, 3\", \"20\u00f73=6, 2\", \"15\u00f73=5, 0\", \"10\u00f78=1, 2\", \"88\u00f78=11, 0\")\n    5  = @(\"34\u00f79=3, 7\", \"14\u00f78=1, 6\", \"48\u00f74=12, 0\", \"40\u00f73=13, 1\", \"31\u00f78=3, 7\")\n    9  = @(\"98\u00f73=32, 2\", \"30\u00f74=7, 2\", \"31\u00f76=5, 1\", \"97\u00f79=10, 7\", \"13\u00f73=4, 1\")\n    13 = @(\"82\u00f74=20, 2\", \"24\u00f74=6, 0\", \"64\u00f72=32, 0\", \"33\u00f78=4, 1\", \"88\u00f75=17, 3\")\n    17 = @(\"72\u00f74=18, 0\", \"15\u00f73=5, 0\", \"67\u00f76=11, 1\", \"37\u00f77=5, 2\", \"60\u00f72=30, 0\")\n}\n\nforeach ($row in $dataRows) {\n    $values = $newValues[$row]\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $t.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
